# This script fills in previously-empty "phone" (column B) cells with the
# placeholder value "N/A". For a handful of rows where phone, status (C)
# and priority (D) were all blank, status also becomes "N/A" and priority
# becomes the numeric value 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only column B (phone) needs to be set to "N/A"
$rowsBOnly = @(
    36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51,
    140, 141, 142, 143, 144, 145, 146, 147, 148, 149, 150, 151, 152, 153, 154, 155,
    601
)

foreach ($r in $rowsBOnly) {
    $ws.Cells.Item($r, 2).Value = "N/A"
}

# Rows where columns B, C and D were all empty: B and C become "N/A" and
# D becomes the number 0
$rowsBCD = @(93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 588, 589, 590, 591)

foreach ($r in $rowsBCD) {
    $ws.Cells.Item($r, 2).Value = "N/A"
    $ws.Cells.Item($r, 3).Value = "N/A"
    $ws.Cells.Item($r, 4).Value = 0
}
